# The edit lives in the primary (default) header of section 1 (header1.xml):
#   1) "***TELJES***" -> "***MISZ***"
#   2) the bold "Indulók:" count "166" -> "0"

$d = $word.ActiveDocument

$section = $d.Sections.First
$header = $section.Headers.Item(1)   # wdHeaderFooterPrimary
$headerRange = $header.Range

# 1) Replace the competition-type marker.
$headerRange.Find.Execute("***TELJES***", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "***MISZ***", 2)

# 2) Replace the "Indulók:" (participants) count.
$headerRange.Find.Execute("166", $true, $true, $false, $false, $false, `
                           $true, 1, $false, "0", 2)
